$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Remove the duplicate "Contact / No display for ContactDetail" row (old row 11).
# This shifts rows 12-22 up by one (new rows 11-21).
$ws1.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now set to "Alvearie Team"
$ws1.Range("B9").Value = "Alvearie Team"

# The remaining "Contact" row (old row 10) becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Case Sensitive value was blank, now set to the text "true" (not boolean TRUE).
# A plain Value="true" assignment auto-coerces to a Boolean cell, so instead
# build the literal text via a formula in a scratch cell, copy its computed
# value across with Paste Special (Values only) which preserves the Text
# type, then clear the scratch cell.
$ws1.Range("Z1").Formula = '="true"'
$ws1.Range("Z1").Copy()
$ws1.Range("B14").PasteSpecial(-4163)
$ws1.Range("Z1").ClearContents()
